# Update "想去人数" (want-to-go count) figures on the 展览 and 全部类型 sheets
# to reflect the latest scrape, per commit "Update gh-pages to output generated at 456a3b4".

$wb = $excel.ActiveWorkbook

$sheet1 = $wb.Worksheets.Item("展览")
$sheet1.Range("F2").Value = 14967
$sheet1.Range("F3").Value = 18870
$sheet1.Range("F14").Value = 140
$sheet1.Range("F15").Value = 214
$sheet1.Range("F22").Value = 7843
$sheet1.Range("F26").Value = 60
$sheet1.Range("F27").Value = 1235
$sheet1.Range("F29").Value = 6019
$sheet1.Range("F35").Value = 5382
$sheet1.Range("F36").Value = 1
$sheet1.Range("F37").Value = 4
$sheet1.Range("F39").Value = 44

$sheet4 = $wb.Worksheets.Item("全部类型")
$sheet4.Range("F2").Value = 14967
$sheet4.Range("F3").Value = 18870
$sheet4.Range("F14").Value = 140
$sheet4.Range("F15").Value = 214
$sheet4.Range("F23").Value = 7843
$sheet4.Range("F27").Value = 60
$sheet4.Range("F28").Value = 1235
$sheet4.Range("F32").Value = 6019
$sheet4.Range("F38").Value = 5382
$sheet4.Range("F39").Value = 1
$sheet4.Range("F40").Value = 4
$sheet4.Range("F42").Value = 44
